$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

$values = @{
    2  = 87.72
    3  = 85.58
    4  = 85.15000000000001
    5  = 86.03
    6  = 82.39
    7  = 89.38
    8  = 97.98999999999999
    9  = 113.2
    10 = 107.7
    11 = 92.91
    12 = 79.56
    13 = 67.45
    14 = 54.6
    15 = 35.3
    16 = 20.96
    17 = 19.7
    18 = 36.78
    19 = 71.87
    20 = 88.27
    21 = 109.29
    22 = 116.83
    23 = 134.94
    24 = 125.71
    25 = 107.86
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 14).Value = $values[$row]
}
